$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 15, shifting existing rows 15-20 down to 16-21.
$ws.Rows.Item(15).Insert()

# Copy formatting/style of the date cell from the row below (now row 16) into new row 15
$ws.Range("D16").Copy() | Out-Null
$ws.Range("D15").PasteSpecial(-4122) | Out-Null
$ws.Application.CutCopyMode = $false

# Populate the newly inserted row 15 with data
$ws.Range("A15").Value = 7
$ws.Range("B15").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C15").Value = "Ñuble"
$ws.Range("D15").Value = 44818
$ws.Range("E15").Value = 16
$ws.Range("F15").Value = 100112044
$ws.Range("G15").Value = "Perejil"
$ws.Range("H15").Value = "Sin especificar"
$ws.Range("I15").Value = "Primera"
$ws.Range("J15").Value = 300
$ws.Range("K15").Value = 800
$ws.Range("L15").Value = 900
$ws.Range("M15").Value = 850
$ws.Range("N15").Value = "$/atado 0,5 a 1 kilo"
$ws.Range("O15").Value = "Región del Maule"
$ws.Range("P15").Value = 850
$ws.Range("Q15").Value = 1
$ws.Range("R15").Value = "Hortaliza"

$wb.Save()
